# Generate Report for Handoff
#
# Refreshes the localization-status report:
#  - bumps the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#    timestamps for the rows whose handoff xliff was just (re)generated
#  - stamps their "Priority" column with "ht" (handoff type) now that the
#    handoff file exists
#
# Affected data rows (1-based sheet rows 8,9,10,12,13,14 -- row 11 is a
# different file whose priority/date aren't touched by this handoff).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$rows = 8,9,10,12,13,14

foreach ($r in $rows) {
    # New handoff xliff generation timestamps, per locale.
    $wsOverview.Range("G" + $r).Value = "2016-08-25 12:22:11"
    $wsZhCn.Range("H" + $r).Value     = "2016-08-25 12:21:59"
    $wsDeDe.Range("H" + $r).Value     = "2016-08-25 12:22:11"

    # Handoff type is now known for these rows.
    $wsZhCn.Range("E" + $r).Value = "ht"
    $wsDeDe.Range("E" + $r).Value = "ht"
}
